# Daily attendance processing - swap the order of "Recorded By" entries
# from "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for every matching cell in column G ("Recorded By").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
